# Apply the "add lower/upper risk, and set ne_path" edit:
#   - Rename the single existing sheet to "line_risk" and replace its
#     (historical/rcp8p5) x (0..4) rows with (historical/rcp8p5) x
#     (mean_risk/lower_risk/upper_risk) summary rows.
#   - Add four more sheets (plant_risk, substation_risk, tower_risk,
#     pole_risk) with the same header/label layout and their own values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper-ish inline logic: style for the header row (B1:C1) and the
# label column (A2:A4) - bold font, thin border all around, centered
# horizontally, top-aligned vertically (mirrors the workbook's existing
# "s=1" cell style).
# ---------------------------------------------------------------------

function Format-RiskLabels {
    param($range)
    $range.Font.Bold = $true
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
}

function Write-RiskSheet {
    param($ws, $meanHist, $meanRcp, $lowerHist, $lowerRcp, $upperHist, $upperRcp)

    $ws.Range("A2").Value = "mean_risk"
    $ws.Range("B2").Value = $meanHist
    $ws.Range("C2").Value = $meanRcp

    $ws.Range("A3").Value = "lower_risk"
    $ws.Range("B3").Value = $lowerHist
    $ws.Range("C3").Value = $lowerRcp

    $ws.Range("A4").Value = "upper_risk"
    $ws.Range("B4").Value = $upperHist
    $ws.Range("C4").Value = $upperRcp

    Format-RiskLabels $ws.Range("A2:A4")
}

# ---------------------------------------------------------------------
# Sheet 1: Sheet1 -> line_risk
# ---------------------------------------------------------------------
$wsLine = $wb.Worksheets.Item(1)
$wsLine.Name = "line_risk"

# Drop the old rows 5 and 6 entirely so the sheet dimension shrinks back
# to A1:C4 (delete the higher-numbered row first so indices stay valid).
$wsLine.Rows.Item(6).Delete()
$wsLine.Rows.Item(5).Delete()

Write-RiskSheet $wsLine 36975.34439514048 54009.1779881274 19522.61815432023 26179.53303909169 67387.80503687025 104458.0386062727

# ---------------------------------------------------------------------
# Sheet 2: plant_risk
# ---------------------------------------------------------------------
$wsPlant = $wb.Worksheets.Add($null, $wsLine)
$wsPlant.Name = "plant_risk"

$wsPlant.Range("B1").Value = "historical"
$wsPlant.Range("C1").Value = "rcp8p5"
Format-RiskLabels $wsPlant.Range("B1:C1")

Write-RiskSheet $wsPlant 642729008.5274855 853671862.0474365 482046756.395614 640253896.5355773 803411260.6593568 1067089827.559296

# ---------------------------------------------------------------------
# Sheet 3: substation_risk
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Add($null, $wsPlant)
$wsSub.Name = "substation_risk"

$wsSub.Range("B1").Value = "historical"
$wsSub.Range("C1").Value = "rcp8p5"
Format-RiskLabels $wsSub.Range("B1:C1")

Write-RiskSheet $wsSub 6895997.694161309 8987723.862688707 5171998.270620981 6740792.897016531 8619997.117701637 11234654.82836089

# ---------------------------------------------------------------------
# Sheet 4: tower_risk
# ---------------------------------------------------------------------
$wsTower = $wb.Worksheets.Add($null, $wsSub)
$wsTower.Name = "tower_risk"

$wsTower.Range("B1").Value = "historical"
$wsTower.Range("C1").Value = "rcp8p5"
Format-RiskLabels $wsTower.Range("B1:C1")

Write-RiskSheet $wsTower 132315.8115881768 179535.9734914216 99236.85869113258 134651.9801185662 165394.764485221 224419.966864277

# ---------------------------------------------------------------------
# Sheet 5: pole_risk
# ---------------------------------------------------------------------
$wsPole = $wb.Worksheets.Add($null, $wsTower)
$wsPole.Name = "pole_risk"

$wsPole.Range("B1").Value = "historical"
$wsPole.Range("C1").Value = "rcp8p5"
Format-RiskLabels $wsPole.Range("B1:C1")

Write-RiskSheet $wsPole 0 0.4835829196380592 0 0.3626871897285443 0 0.604478649547574

# Select line_risk as the active sheet (matches activeTab="0" in the target).
$wsLine.Activate()
